# TLWP-845 - Add columns to pipeline report
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Opportunities with providers"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Opportunities with providers")

# New header cells (G1:L1) - copy the header formatting from an existing
# header cell first so the new cells pick up the bold/left-aligned style,
# then overwrite the value (PasteSpecial formats only, Value stays intact
# afterwards).
$headerSource = $ws1.Range("F1")
$newHeaders = @(
    @{ Cell = "G1"; Text = "Provider contact name" },
    @{ Cell = "H1"; Text = "Provider contact email" },
    @{ Cell = "I1"; Text = "Provider contact telephone" },
    @{ Cell = "J1"; Text = "Secondary contact name" },
    @{ Cell = "K1"; Text = "Secondary contact email" },
    @{ Cell = "L1"; Text = "Secondary contact telephone" }
)

foreach ($h in $newHeaders) {
    $headerSource.Copy() | Out-Null
    $ws1.Range($h.Cell).PasteSpecial(-4122) | Out-Null
    $ws1.Range($h.Cell).Value = $h.Text
}

# Column widths (characters). Columns A-E keep their existing (auto-fit)
# widths untouched; F and the new G:L columns get an explicit width closest
# to the template.
$ws1.Range("F1").ColumnWidth = 24.333333333333332
$ws1.Range("G1").ColumnWidth = 21.833333333333332
$ws1.Range("H1").ColumnWidth = 21.833333333333332
$ws1.Range("I1").ColumnWidth = 26.5
$ws1.Range("J1").ColumnWidth = 24.5
$ws1.Range("K1").ColumnWidth = 24.166666666666668
$ws1.Range("L1").ColumnWidth = 28.833333333333332

# Row height for the header row.
$ws1.Rows.Item(1).RowHeight = 15.75

# Selection / scroll position, matching the saved view in the template.
$ws1.Range("L1").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2: "Opportunities with no providers"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Opportunities with no providers")

# Insert a new "Job role" column before the existing "Number of students
# wanted" column, shifting "Reason no providers chosen" to column D.
$ws2.Range("B1").EntireColumn.Insert(-4161) | Out-Null

$b1Source = $ws2.Range("A1")
$b1Source.Copy() | Out-Null
$ws2.Range("B1").PasteSpecial(-4122) | Out-Null
$ws2.Range("B1").Value = "Job role"

# Only the brand new column needs an explicit width; C and D keep the
# auto-fit widths (and bestFit flag) they already had as B and C.
$ws2.Range("B1").ColumnWidth = 29.5

$ws2.Range("B15").Select() | Out-Null

# ---------------------------------------------------------------------------
# Re-activate sheet 1 so it remains the selected/visible tab on save, as in
# the template (tabSelected="1" on sheet 1 only).
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("L1").Select() | Out-Null
